$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.836.16'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '2.282.67'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'505.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.95%  '
$ws.Range("D6").Value = "'128.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("D9").Value = '2.297.15'
$ws.Range("D10").Value = "'0.0967"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("E12").Value = '  +4.55%  '
$ws.Range("E13").Value = '  +5.31%  '
$ws.Range("D14").Value = "'23.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.33%  '
$ws.Range("D15").Value = '2.685.38'
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '54.875.41'
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("D18").Value = '2.270.32'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("D19").Value = "'10.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.25%  '
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("D21").Value = "'307.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.56%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = "'60.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.90%  '
$ws.Range("D25").Value = "'0.996"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("E26").Value = '  +0.86%  '
$ws.Range("D27").Value = "'7.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.69%  '
$ws.Range("D28").Value = "'171.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = "'6.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.20%  '
$ws.Range("D30").Value = '0.0₃0703'
$ws.Range("E30").Value = '  +3.40%  '
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("E32").Value = '  +2.94%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  +1.23%  '
$ws.Range("D35").Value = "'0.993"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.33%  '
$ws.Range("D36").Value = "'0.918"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.06%  '
$ws.Range("D37").Value = "'1.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.36%  '
$ws.Range("E38").Value = '  +2.34%  '
$ws.Range("D39").Value = "'36.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.19%  '
$ws.Range("D40").Value = "'0.376"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.22%  '
$ws.Range("D41").Value = "'1.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("E42").Value = '  +6.05%  '
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("D44").Value = "'126.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("E45").Value = '  +2.21%  '
$ws.Range("D46").Value = "'249.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.11%  '
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("D49").Value = "'0.376"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.36%  '
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("D51").Value = "'10.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.51%  '
